# Commit: "Completed Up to - Results section on report"
# The underlying edit is a spelling fix: "saloon" -> "salon" throughout the
# chat-text entries stored in column A of the only worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A1:A481")
$rng.Replace("saloon", "salon") | Out-Null

# Reflect the cursor/scroll position recorded in the saved workbook after
# the edit was made.
$excel.ActiveWindow.ScrollRow = 154
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G172").Select() | Out-Null
